$d = $word.ActiveDocument

# Locate the two anchor paragraphs by their placeholder text:
#  - $afterPara : the "Liebe ... [Adressat.Vorname]" paragraph, after which the
#                 new [Adressat.Anschrift] line must be inserted.
#  - $templatePara : the last "[Fall.Versicherungsnummer] ..." list paragraph,
#                 whose paragraph (pPr) and run (rPr) formatting exactly match
#                 what the new paragraph needs.
$afterPara = $null
$templatePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Liebe*Adressat.Vorname]*") {
        $afterPara = $p
    }
    if ($t -like "*Fall.Versicherungsnummer*") {
        $templatePara = $p
    }
}

# Build a new paragraph that inherits the template paragraph's formatting by
# inserting it directly in front of the template paragraph. After the insert,
# $templatePara itself refers to the freshly-created (still empty) paragraph,
# while the original content shifted down to the next paragraph.
$templatePara.Range.InsertParagraphBefore()
$newPara = $templatePara
$newPara.Range.Text = "[Adressat.Anschrift]"

# Move the freshly formatted paragraph to its proper place, right after the
# "Liebe ... [Adressat.Vorname]" paragraph.
$newPara.Range.Cut()
$dest = $d.Range($afterPara.Range.End, $afterPara.Range.End)
$dest.Paste()

Write-Host "done"
